$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update "last updated" timestamp (row 1) ---
$ws.Range("A1").Value = "Datos actualizados a 20 de Abril de 2020 a las 22:22"

# --- Row 4: Estados Unidos ---
$ws.Range("B4").Value = 784560
$ws.Range("C4").Value = 19924
$ws.Range("E4").Value = 670946
$ws.Range("G4").Value = 1269
$ws.Range("H4").Value = 41844

# --- Row 8: Alemania ---
$ws.Range("B8").Value = 146653
$ws.Range("C8").Value = 911
$ws.Range("E8").Value = 50447

# --- Rows 14 & 15: Brasil overtakes Belgica in the ranking ---
# Row 14 becomes Brasil (new, higher numbers)
$ws.Range("A14").Value = "Brasil"
$ws.Range("B14").Value = 40581
$ws.Range("C14").Value = 1927
$ws.Range("D14").Value = 22130
$ws.Range("E14").Value = 15606
$ws.Range("F14").Value = 7919
$ws.Range("G14").Value = 383
$ws.Range("H14").Value = 2845

# Row 15 becomes Belgica (its previous, unchanged figures)
$ws.Range("A15").Value = "Belgica"
$ws.Range("B15").Value = 39983
$ws.Range("C15").Value = 1487
$ws.Range("D15").Value = 8895
$ws.Range("E15").Value = 25260
$ws.Range("F15").Value = 1071
$ws.Range("G15").Value = 145
$ws.Range("H15").Value = 5828

# --- Row 86: Tunez ---
$ws.Range("B86").Value = 884
$ws.Range("C86").Value = 5
$ws.Range("E86").Value = 698
$ws.Range("F86").Value = 34

# --- Row 91: Principado de Andorra ---
$ws.Range("E91").Value = 432
$ws.Range("G91").Value = 1
$ws.Range("H91").Value = 37

# --- Rows 151 & 152: Maldivas overtakes Cabo Verde in the ranking ---
# Row 151 becomes Maldivas (new, higher numbers)
$ws.Range("A151").Value = "Maldivas"
$ws.Range("B151").Value = 69
$ws.Range("C151").Value = 17
$ws.Range("D151").Value = 16
$ws.Range("E151").Value = 53
$ws.Range("F151").Value = 1
$ws.Range("G151").Value = 0
$ws.Range("H151").Value = 0

# Row 152 becomes Cabo Verde (its previous, unchanged figures)
$ws.Range("A152").Value = "Cabo Verde"
$ws.Range("B152").Value = 67
$ws.Range("C152").Value = 6
$ws.Range("D152").Value = 1
$ws.Range("E152").Value = 65
$ws.Range("F152").Value = 0
$ws.Range("G152").Value = 0
$ws.Range("H152").Value = 1
